$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 38.33987524840568
$ws.Range("C2").Value = 39.22271935221791
$ws.Range("D2").Value = 37.43837867784747
$ws.Range("E2").Value = 38.34010229605165
$ws.Range("F2").Value = 38.44202087206734
$ws.Range("G2").Value = 37.27085455945915
$ws.Range("H2").Value = 41.12918685901104
$ws.Range("I2").Value = 34.4046478788235
$ws.Range("J2").Value = 38.66587587865975
$ws.Range("K2").Value = 38.41084676549933
$ws.Range("L2").Value = 38.348474045005
$ws.Range("M2").Value = 38.74644130675284
$ws.Range("N2").Value = 22.52128341190802
$ws.Range("O2").Value = 33.85002466909647
$ws.Range("P2").Value = 41.79484410913059
$ws.Range("Q2").Value = 34.45218478079896
